$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (09-nov) before column DL ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Inserting the column shifts DL:EP (and their data) one column to the right,
# to DM:EQ, exactly like picking "Insert" on the column header in Excel.
$ws1.Columns("DL").Insert()

# Header cell for the freshly inserted column.
$ws1.Range("DL1").Value = "09-nov"

# The data rows (2-25) for that brand-new column have no observations yet,
# matching the "-" placeholder used by the neighbouring not-yet-available
# date columns (DF:DK).
$ws1.Range("DL2:DL25").Value = "-"

# --- Sheet "Gaz": append the new daily price row ---
$ws2 = $wb.Worksheets.Item("Gaz")
$a2 = $ws2.Cells.Item(145, 1)
# Assigning a literal "2025-11-07" string directly would be auto-converted to
# a date serial by the smart-typing layer; going through a text formula and
# then "paste values" keeps it as a genuine text cell (same as every other
# date cell already in this column), without leaving a quotePrefix/NumberFormat
# style behind.
$a2.Formula = "=""2025-11-07"""
$a2.Copy()
$a2.PasteSpecial(-4163)
$ws2.Application.CutCopyMode = $false
$ws2.Cells.Item(145, 2).Value = 29.74

# --- Sheet "CO2": append the new daily price row ---
$ws3 = $wb.Worksheets.Item("CO2")
$a3 = $ws3.Cells.Item(145, 1)
$a3.Formula = "=""2025-11-07"""
$a3.Copy()
$a3.PasteSpecial(-4163)
$ws3.Application.CutCopyMode = $false
$ws3.Cells.Item(145, 2).Value = 79.36
